# Auto-generated edit script: updates FFXIV leve-profit market-data sheets
# (currentAveragePrice / LevePrice / LeveProfit columns H-N) per scheduled
# runner refresh. Values only; no structural/formatting changes.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 430
$ws.Range("I6").Value = 145
$ws.Range("K6").Value = 435
$ws.Range("M6").Value = -323
$ws.Range("H32").Value = 700
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 700
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 700
$ws.Range("M32").Value = ""
$ws.Range("N32").Value = -1352
$ws.Range("H43").Value = 100
$ws.Range("I43").Value = 100
$ws.Range("K43").Value = 100
$ws.Range("M43").Value = -31
$ws.Range("H55").Value = 266.23077
$ws.Range("I55").Value = 596.6667
$ws.Range("J55").Value = 91.29412000000001
$ws.Range("K55").Value = 596.6667
$ws.Range("L55").Value = 91.29412000000001
$ws.Range("M55").Value = -382.6667
$ws.Range("N55").Value = -519.29412
$ws.Range("H58").Value = 4024.5
$ws.Range("J58").Value = 7672.3335
$ws.Range("L58").Value = 23017.0005
$ws.Range("N58").Value = -23317.0005
$ws.Range("H62").Value = 2780.318
$ws.Range("I62").Value = 2016.5385
$ws.Range("J62").Value = 3883.5557
$ws.Range("K62").Value = 2016.5385
$ws.Range("L62").Value = 3883.5557
$ws.Range("M62").Value = -1392.5385
$ws.Range("N62").Value = -5131.5557
$ws.Range("H65").Value = 2780.318
$ws.Range("I65").Value = 2016.5385
$ws.Range("J65").Value = 3883.5557
$ws.Range("K65").Value = 10082.6925
$ws.Range("L65").Value = 19417.7785
$ws.Range("M65").Value = -6962.692500000001
$ws.Range("N65").Value = -25657.7785
$ws.Range("H76").Value = 3542.2727
$ws.Range("I76").Value = 3327.5
$ws.Range("K76").Value = 3327.5
$ws.Range("M76").Value = -3012.5
$ws.Range("H79").Value = 3542.2727
$ws.Range("I79").Value = 3327.5
$ws.Range("K79").Value = 3327.5
$ws.Range("M79").Value = -2235.5
$ws.Range("H112").Value = 1148.4615
$ws.Range("J112").Value = 1148.4615
$ws.Range("L112").Value = 3445.3845
$ws.Range("N112").Value = -5661.3845
$ws.Range("H129").Value = 667284.6
$ws.Range("J129").Value = 1000797.7
$ws.Range("L129").Value = 3002393.1
$ws.Range("N129").Value = -3012393.1
$ws.Range("H131").Value = 2161.25
$ws.Range("I131").Value = 322.5
$ws.Range("K131").Value = 967.5
$ws.Range("M131").Value = 4072.5
$ws.Range("H132").Value = 4352.647
$ws.Range("I132").Value = 4562.1875
$ws.Range("J132").Value = 1000
$ws.Range("K132").Value = 13686.5625
$ws.Range("L132").Value = 3000
$ws.Range("M132").Value = -11156.5625
$ws.Range("N132").Value = -8060
$ws.Range("H137").Value = 1655.2
$ws.Range("I137").Value = 1398.421
$ws.Range("K137").Value = 4195.263
$ws.Range("M137").Value = -1645.263
$ws.Range("H138").Value = 2374.791
$ws.Range("I138").Value = 2735.5454
$ws.Range("J138").Value = 2303.9285
$ws.Range("K138").Value = 8206.636200000001
$ws.Range("L138").Value = 6911.7855
$ws.Range("M138").Value = -3066.636200000001
$ws.Range("N138").Value = -17191.7855

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4313.8433
$ws.Range("I32").Value = 4387.8774
$ws.Range("J32").Value = 2500
$ws.Range("K32").Value = 4387.8774
$ws.Range("L32").Value = 2500
$ws.Range("M32").Value = -4100.8774
$ws.Range("N32").Value = -3074
$ws.Range("H45").Value = 1942.1562
$ws.Range("I45").Value = 2220.125
$ws.Range("J45").Value = 1664.1875
$ws.Range("K45").Value = 2220.125
$ws.Range("L45").Value = 1664.1875
$ws.Range("M45").Value = -1843.125
$ws.Range("N45").Value = -2418.1875
$ws.Range("H61").Value = 1743.697
$ws.Range("I61").Value = 1512.9656
$ws.Range("K61").Value = 1512.9656
$ws.Range("M61").Value = -1300.9656
$ws.Range("H132").Value = 20709.482
$ws.Range("I132").Value = 2139.8948
$ws.Range("K132").Value = 6419.6844
$ws.Range("M132").Value = -3889.6844
$ws.Range("H136").Value = 1743.697
$ws.Range("I136").Value = 1512.9656
$ws.Range("K136").Value = 4538.8968
$ws.Range("M136").Value = -1988.8968

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 943.2778
$ws.Range("I20").Value = 940.9167
$ws.Range("K20").Value = 940.9167
$ws.Range("M20").Value = -693.9167
$ws.Range("H134").Value = 3730.4482
$ws.Range("I134").Value = 4115.32
$ws.Range("J134").Value = 1325
$ws.Range("K134").Value = 12345.96
$ws.Range("L134").Value = 3975
$ws.Range("M134").Value = -9810.959999999999
$ws.Range("N134").Value = -9045

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 19499.5
$ws.Range("J28").Value = 19499.5
$ws.Range("L28").Value = 19499.5
$ws.Range("N28").Value = -19989.5
$ws.Range("H105").Value = 3470.3333
$ws.Range("I105").Value = 400
$ws.Range("J105").Value = 5005.5
$ws.Range("K105").Value = 400
$ws.Range("L105").Value = 5005.5
$ws.Range("M105").Value = 1347
$ws.Range("N105").Value = -8499.5
$ws.Range("H107").Value = 869
$ws.Range("I107").Value = 251.66667
$ws.Range("J107").Value = 1486.3334
$ws.Range("K107").Value = 251.66667
$ws.Range("L107").Value = 1486.3334
$ws.Range("M107").Value = 1668.33333
$ws.Range("N107").Value = -5326.3334
$ws.Range("H134").Value = 1251.0834
$ws.Range("I134").Value = 1227.1818
$ws.Range("J134").Value = 1514
$ws.Range("K134").Value = 3681.5454
$ws.Range("L134").Value = 4542
$ws.Range("M134").Value = -1146.5454
$ws.Range("N134").Value = -9612

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 36.76923
$ws.Range("I2").Value = 27.285715
$ws.Range("J2").Value = 47.833332
$ws.Range("K2").Value = 163.71429
$ws.Range("L2").Value = 286.999992
$ws.Range("M2").Value = -50.71429000000001
$ws.Range("N2").Value = -512.999992
$ws.Range("H7").Value = 45
$ws.Range("I7").Value = 40
$ws.Range("J7").Value = 50
$ws.Range("K7").Value = 120
$ws.Range("L7").Value = 150
$ws.Range("M7").Value = -8
$ws.Range("N7").Value = -374
$ws.Range("H22").Value = 50550
$ws.Range("I22").Value = 50550
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 151650
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -151481
$ws.Range("N22").Value = ""
$ws.Range("H27").Value = 50550
$ws.Range("I27").Value = 50550
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 151650
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = -151548
$ws.Range("N27").Value = ""
$ws.Range("H38").Value = 128.2
$ws.Range("I38").Value = 110.25
$ws.Range("J38").Value = 200
$ws.Range("K38").Value = 330.75
$ws.Range("L38").Value = 600
$ws.Range("M38").Value = 16.25
$ws.Range("N38").Value = -1294
$ws.Range("I92").Value = 31250676
$ws.Range("J92").Value = 600
$ws.Range("K92").Value = 93752028
$ws.Range("L92").Value = 1800
$ws.Range("M92").Value = -93750780
$ws.Range("N92").Value = -4296
$ws.Range("H113").Value = 756
$ws.Range("J113").Value = 761.8
$ws.Range("L113").Value = 2285.4
$ws.Range("N113").Value = -6625.4
$ws.Range("H131").Value = 747.39
$ws.Range("J131").Value = 747.4141
$ws.Range("L131").Value = 2242.2423
$ws.Range("N131").Value = -12322.2423
$ws.Range("H132").Value = 1090.5
$ws.Range("J132").Value = 1297.1666
$ws.Range("L132").Value = 11674.4994
$ws.Range("N132").Value = -16734.4994

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 645.3333
$ws.Range("I107").Value = 234.1
$ws.Range("J107").Value = 2701.5
$ws.Range("K107").Value = 234.1
$ws.Range("L107").Value = 2701.5
$ws.Range("M107").Value = 1685.9
$ws.Range("N107").Value = -6541.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 929.65
$ws.Range("I46").Value = 924.4167
$ws.Range("J46").Value = 937.5
$ws.Range("K46").Value = 924.4167
$ws.Range("L46").Value = 937.5
$ws.Range("M46").Value = -736.4167
$ws.Range("N46").Value = -1313.5
$ws.Range("H97").Value = 17172
$ws.Range("J97").Value = 17172
$ws.Range("L97").Value = 17172
$ws.Range("N97").Value = -19154
$ws.Range("H136").Value = 1004.41174
$ws.Range("I136").Value = 953.9091
$ws.Range("J136").Value = 1097
$ws.Range("K136").Value = 2861.7273
$ws.Range("L136").Value = 3291
$ws.Range("M136").Value = -311.7273
$ws.Range("N136").Value = -8391

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 939.1667
$ws.Range("I113").Value = 1477.8
$ws.Range("J113").Value = 265.875
$ws.Range("K113").Value = 4433.4
$ws.Range("L113").Value = 797.625
$ws.Range("M113").Value = -2263.4
$ws.Range("N113").Value = -5137.625

